$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "1.00" or "598.92" would otherwise be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.927.08"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "3.825.30"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "598.92"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "178.50"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "3.827.00"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  -5.14%  "
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("D13").Value = "38.58"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "4.460.83"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "3.825.05"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "67.920.77"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D20").Value = "16.47"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "491.73"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D24").Value = "0.0000151"
$ws.Range("E24").Value = "  +10.49%  "
$ws.Range("D25").Value = "84.87"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  -6.61%  "
$ws.Range("D31").Value = "2.45"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "33.09"
$ws.Range("E32").Value = "  +8.30%  "
$ws.Range("D33").Value = "7.77"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Value = "0.137"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "5.80"
$ws.Range("E38").Value = "  -5.80%  "
$ws.Range("D40").Value = "455.54"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D45").Value = "41.45"
$ws.Range("E45").Value = "  -8.22%  "
$ws.Range("D46").Value = "2.852.54"
$ws.Range("E46").Value = "  -3.53%  "
$ws.Range("D47").Value = "140.88"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D50").Value = "26.24"
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("D51").Value = "23.78"
$ws.Range("E51").Value = "  +8.71%  "

$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E26").Value = "  -8.11%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("E49").Value = "  -2.86%  "

# Rows 41-43: source data reshuffled (OKB / Stacks / dogwifhat rotated down one slot)
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "49.28"
$ws.Range("E41").Value = "  -0.85%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.02"
$ws.Range("E42").Value = "  -2.61%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  -5.04%  "
